# Commit: "Added parsing for holes"
#
# The course name used for the second 9/18-hole course was stored with a
# typo ("discgoldbana" instead of "discgolfbana"). Fix the name for every
# hole-row belonging to that course (rows 11-28, column A / "name").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("A11:A28").Value = "almhults discgolfbana"

# Reflect where the user ended up working in the sheet (bottom of the
# data, near the rows that were just edited).
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("G27").Select()
